$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column cells originally hold inline-string (text) values such as "62.26".
# Setting .Value directly with a numeric-looking string would make Excel convert the cell
# to a real number, changing its stored representation. Temporarily force text format on
# each Price cell we touch, then restore the default style afterward.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "38.784.13"
$ws.Range("E2").Value = "  +0.32%  "

$ws.Range("D3").Value = "2.104.91"
$ws.Range("E3").Value = "  +0.37%  "

$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("E5").Value = "  -0.57%  "

$ws.Range("D7").Value = "62.22"
$ws.Range("E7").Value = "  +0.99%  "

$ws.Range("E9").Value = "  +2.00%  "

$ws.Range("E10").Value = "  -0.44%  "

$ws.Range("E11").Value = "  -0.95%  "

$ws.Range("D12").Value = "15.76"
$ws.Range("E12").Value = "  +6.25%  "

$ws.Range("D13").Value = "2.418.29"
$ws.Range("E13").Value = "  +0.40%  "

$ws.Range("D14").Value = "22.09"
$ws.Range("E14").Value = "  -1.48%  "

$ws.Range("D15").Value = "0.809"
$ws.Range("E15").Value = "  +3.12%  "

$ws.Range("D16").Value = "5.53"
$ws.Range("E16").Value = "  +1.01%  "

$ws.Range("D17").Value = "2.103.95"
$ws.Range("E17").Value = "  +0.18%  "

$ws.Range("D18").Value = "38.794.24"
$ws.Range("E18").Value = "  +0.53%  "

$ws.Range("D19").Value = "71.98"
$ws.Range("E19").Value = "  +1.17%  "

$ws.Range("E20").Value = "  +0.94%  "

$ws.Range("E21").Value = "  +0.55%  "

$ws.Range("D22").Value = "228.22"
$ws.Range("E22").Value = "  +0.77%  "

$ws.Range("E23").Value = "  +0.02%  "

$ws.Range("E24").Value = "  -3.27%  "

$ws.Range("E25").Value = "  -0.45%  "

$ws.Range("D26").Value = "9.66"
$ws.Range("E26").Value = "  +2.06%  "

$ws.Range("D27").Value = "172.35"
$ws.Range("E27").Value = "  +1.08%  "

$ws.Range("D28").Value = "0.139"
$ws.Range("E28").Value = "  +5.33%  "

$ws.Range("E29").Value = "  +4.54%  "

$ws.Range("E30").Value = "  +1.11%  "

$ws.Range("D31").Value = "2.52"
$ws.Range("E31").Value = "  +10.31%  "

$ws.Range("E32").Value = "  +0.54%  "

$ws.Range("E33").Value = "  +1.39%  "

$ws.Range("E34").Value = "  -0.06%  "

$ws.Range("D35").Value = "7.04"
$ws.Range("E35").Value = "  +6.99%  "

$ws.Range("D36").Value = "0.0621"
$ws.Range("E36").Value = "  +2.15%  "

$ws.Range("E37").Value = "  +0.27%  "

$ws.Range("E38").Value = "  +1.13%  "

$ws.Range("E39").Value = "  +0.06%  "

$ws.Range("E40").Value = "  -3.60%  "

$ws.Range("D41").Value = "102.75"
$ws.Range("E41").Value = "  +2.61%  "

$ws.Range("E42").Value = "  +3.63%  "

$ws.Range("D43").Value = "1.528.99"
$ws.Range("E43").Value = "  -1.13%  "

$ws.Range("D44").Value = "1.20"
$ws.Range("E44").Value = "  +6.77%  "

$ws.Range("E45").Value = "  -1.09%  "

$ws.Range("D46").Value = "7.75"
$ws.Range("E46").Value = "  +1.23%  "

$ws.Range("D47").Value = "0.0912"
$ws.Range("E47").Value = "  -0.42%  "

$ws.Range("E48").Value = "  -0.90%  "

$ws.Range("D49").Value = "1.06"

$ws.Range("E50").Value = "  -0.39%  "

$ws.Range("D51").Value = "2.304.43"
$ws.Range("E51").Value = "  +0.44%  "

# Restore the default (Normal) cell style on the Price cells we touched.
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Style = "Normal"
$ws.Range("D7").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("D51").Style = "Normal"
